$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.895.79'
$ws.Range("E2").Value = '  -0.22%  '

$ws.Range("D3").Value = '1.634.54'
$ws.Range("E3").Value = '  -0.41%  '

$ws.Range("D4").Formula = "'1.002"
$ws.Range("E4").Value = '  -0.33%  '

$ws.Range("D5").Formula = "'214.24"
$ws.Range("E5").Value = '  -0.53%  '

$ws.Range("D6").Formula = "'0.5054"
$ws.Range("E6").Value = '  -0.52%  '

$ws.Range("E7").Value = '  -0.25%  '

$ws.Range("D8").Formula = "'0.2571"
$ws.Range("E8").Value = '  +0.42%  '

$ws.Range("D9").Formula = "'0.06350"
$ws.Range("E9").Value = '  -0.46%  '

$ws.Range("D10").Formula = "'19.68"
$ws.Range("E10").Value = '  +1.14%  '

$ws.Range("D11").Formula = "'0.07738"
$ws.Range("E11").Value = '  -0.53%  '

$ws.Range("D12").Formula = "'4.281"
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").Value = '1.636.12'
$ws.Range("E13").Value = '  -0.39%  '

$ws.Range("D14").Formula = "'0.5440"
$ws.Range("E14").Value = '  -0.37%  '

$ws.Range("D15").Value = '0.0₅7727'

$ws.Range("D16").Formula = "'64.03"
$ws.Range("E16").Value = '  -0.56%  '

$ws.Range("D17").Value = '25.915.36'
$ws.Range("E17").Value = '  -0.37%  '

$ws.Range("D18").Formula = "'1.002"

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Formula = "'195.34"
$ws.Range("E19").Value = '  -1.38%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Formula = "'4.431"
$ws.Range("E20").Value = '  -0.20%  '

$ws.Range("D21").Formula = "'9.911"
$ws.Range("E21").Value = '  -0.50%  '

$ws.Range("D22").Formula = "'6.119"
$ws.Range("E22").Value = '  +1.11%  '

$ws.Range("D23").Formula = "'1.003"
$ws.Range("E23").Value = '  -0.31%  '

$ws.Range("D24").Formula = "'1.889"
$ws.Range("E24").Value = '  +0.66%  '

$ws.Range("D25").Formula = "'142.92"
$ws.Range("E25").Value = '  +1.43%  '

$ws.Range("D26").Formula = "'0.1242"
$ws.Range("E26").Value = '  +8.27%  '

$ws.Range("D27").Formula = "'6.823"
$ws.Range("E27").Value = '  -0.84%  '

$ws.Range("D28").Formula = "'15.62"
$ws.Range("E28").Value = '  -0.74%  '

$ws.Range("D29").Formula = "'1.236"
$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("D30").Formula = "'0.04867"
$ws.Range("E30").Value = '  -3.13%  '

$ws.Range("D31").Formula = "'3.236"
$ws.Range("E31").Value = '  -0.69%  '

$ws.Range("D32").Formula = "'3.194"
$ws.Range("E32").Value = '  +0.22%  '

$ws.Range("D33").Formula = "'1.546"
$ws.Range("E33").Value = '  +0.30%  '

$ws.Range("D34").Formula = "'2.371"
$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("D35").Formula = "'0.9108"
$ws.Range("E35").Value = '  +1.33%  '

$ws.Range("E36").Value = '  -1.02%  '

$ws.Range("D37").Formula = "'0.5502"
$ws.Range("E37").Value = '  +0.13%  '

$ws.Range("D38").Value = '1.122.57'
$ws.Range("E38").Value = '  -0.81%  '

$ws.Range("D39").Formula = "'0.01560"
$ws.Range("E39").Value = '  +0.16%  '

$ws.Range("D40").Formula = "'1.001"
$ws.Range("E40").Value = '  -0.33%  '

$ws.Range("D41").Formula = "'5.592"
$ws.Range("E41").Value = '  -0.44%  '

$ws.Range("D42").Formula = "'0.8038"
$ws.Range("E42").Value = '  -1.92%  '

$ws.Range("D43").Formula = "'98.52"
$ws.Range("E43").Value = '  -1.64%  '

$ws.Range("E44").Value = '  -8.38%  '

$ws.Range("D45").Value = '1.769.59'
$ws.Range("E45").Value = '  -0.53%  '

$ws.Range("D46").Formula = "'0.4481"
$ws.Range("E46").Value = '  -1.13%  '

$ws.Range("D47").Formula = "'1.003"
$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("D48").Formula = "'54.93"
$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("D49").Formula = "'0.05174"
$ws.Range("E49").Value = '  +2.03%  '

$ws.Range("D50").Formula = "'7.494"
$ws.Range("E50").Value = '  +1.45%  '

$ws.Range("E51").Value = '  -0.39%  '
